# SALT_properties.xlsx update — "updated some work 3.16"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "KCl" — fill in rows 9, 11, 12 (extend existing D8:G10 table to
# D8:G12) and rows 15-19 (second small table gets E/F/G data + formulas).
# ---------------------------------------------------------------------
$kcl = $wb.Worksheets.Item("KCl")

$kcl.Range("E9").Value = -631.34152687333403
$kcl.Range("F9").Value = 2.6969799999999999

$kcl.Range("E11").Value = -626.46267805333298
$kcl.Range("F11").Value = -0.80862000000000001
$kcl.Range("G11").Formula = "=G10*((D11/D10)^(3))"

$kcl.Range("E12").Value = -623.38357945333303
$kcl.Range("F12").Value = -1.81626
$kcl.Range("G12").Formula = "=G11*((D12/D11)^(3))"

$kcl.Range("E15").Value = -635.36910887999898
$kcl.Range("F15").Value = 12.896513333333299
$kcl.Range("G15").Value = 6859.28

$kcl.Range("E16").Value = -635.38609069999995
$kcl.Range("F16").Value = 6.2037000000000004
$kcl.Range("G16").Formula = "=G15*((D16/D15)^(3))"

$kcl.Range("G17").Formula = "=G16*((D17/D16)^(3))"
$kcl.Range("G18").Formula = "=G17*((D18/D17)^(3))"
$kcl.Range("G19").Formula = "=G18*((D19/D18)^(3))"

# Chart 1 (scatter of KCl!F/G) — extend the x-axis scale and resize the
# chart frame to match the larger table.
$kclChart = $kcl.ChartObjects(1).Chart
$axis = $kclChart.Axes(1)
$axis.MinimumScale = 6000
$axis.MaximumScale = 10000
$kcl.ChartObjects(1).Height = 288

$kcl.Range("B3:G12").Select()

# ---------------------------------------------------------------------
# Sheet "LiCl_80KCl" — previously empty, now holds a full equilibration
# table (title + D/E/F/G columns), mirroring the KCl sheet's layout.
# ---------------------------------------------------------------------
$lic = $wb.Worksheets.Item("LiCl_80KCl")

$lic.Range("B2").Value = "LiCl-80KCl equilibrations"

$lic.Range("D4").Value = 1000

$lic.Range("D5").Value = 0.95
$lic.Range("E5").Value = -645.49680606666595
$lic.Range("F5").Value = 35.19162
$lic.Range("G5").Value = 5588.41

$lic.Range("D6").Value = 0.97499999999999998
$lic.Range("G6").Formula = "=G5*((D6/D5)^(3))"

$lic.Range("D7").Value = 1
$lic.Range("E7").Value = -648.61861583333302
$lic.Range("F7").Value = 14.8488133333333
$lic.Range("G7").Formula = "=G6*((D7/D6)^(3))"

$lic.Range("D8").Value = 1.0249999999999999
$lic.Range("E8").Value = -648.88683186666697
$lic.Range("F8").Value = 8.2293666666666692
$lic.Range("G8").Formula = "=G7*((D8/D7)^(3))"

$lic.Range("D9").Value = 1.05
$lic.Range("E9").Value = -646.176745193334
$lic.Range("F9").Value = 4.9316599999999999
$lic.Range("G9").Formula = "=G8*((D9/D8)^(3))"

$lic.Range("D10").Value = 1.075
$lic.Range("G10").Formula = "=G9*((D10/D9)^(3))"

$lic.Range("D11").Value = 1.1000000000000001
$lic.Range("G11").Formula = "=G10*((D11/D10)^(3))"

$lic.Range("D12").Value = 1.125
$lic.Range("G12").Formula = "=G11*((D12/D11)^(3))"

$lic.Range("D13").Value = 1.1499999999999999
$lic.Range("G13").Formula = "=G12*((D13/D12)^(3))"

$lic.Range("E9:F9").Select()

# ---------------------------------------------------------------------
# Sheet "NaCl_MgCl2" — append a second equilibration block in rows 12-18.
# ---------------------------------------------------------------------
$nacl = $wb.Worksheets.Item("NaCl_MgCl2")

$nacl.Range("B12").Value = 900

$nacl.Range("C13").Value = "a0"
$nacl.Range("D13").Value = "E"
$nacl.Range("E13").Value = "P"
$nacl.Range("F13").Value = "V"

$nacl.Range("C14").Value = 0.9
$nacl.Range("F14").Formula = "=F15*(C14/C15)^3"

$nacl.Range("C15").Value = 0.95
$nacl.Range("F15").Formula = "=F16*(C15/C16)^3"

$nacl.Range("C16").Value = 1
$nacl.Range("F16").Formula = "=F17*(C16/C17)^3"

$nacl.Range("C17").Value = 1.05
$nacl.Range("F17").Formula = "=F18*(C17/C18)^3"

$nacl.Range("C18").Value = 1.1000000000000001
$nacl.Range("D18").Value = -639.56252220666704
$nacl.Range("E18").Value = -2.83110666666667
$nacl.Range("F18").Value = 10144.73

$nacl.Range("G27").Select()

# Active sheet stays NaCl_MgCl2 (already the selected tab).
$nacl.Activate()
